$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 581, shifting existing rows (old 581..673) down to (582..674)
$ws.Rows.Item(581).Insert()

# Populate the newly inserted row 581 with the new record
$ws.Cells.Item(581,1).Value = 3
$ws.Cells.Item(581,2).Value = "Femacal de La Calera"
$ws.Cells.Item(581,3).Value = "Coquimbo"
$ws.Cells.Item(581,4).Value = 45180
$ws.Cells.Item(581,5).Value = 5
$ws.Cells.Item(581,6).Value = 100112031
$ws.Cells.Item(581,7).Value = "Poroto verde"
$ws.Cells.Item(581,8).Value = "Sin especificar"
$ws.Cells.Item(581,9).Value = "Primera"
$ws.Cells.Item(581,10).Value = 78
$ws.Cells.Item(581,11).Value = 27000
$ws.Cells.Item(581,12).Value = 28000
$ws.Cells.Item(581,13).Value = 27487
$ws.Cells.Item(581,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(581,15).Value = "Perú"
$ws.Cells.Item(581,16).Value = 1099
$ws.Cells.Item(581,17).Value = 25
$ws.Cells.Item(581,18).Value = "Hortaliza"
